$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set in_service (column E) to TRUE for rows 10 through 15 (extr3..extr8)
$ws.Range("E10:E15").Value = $true
